$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sits right after "Sample
#    Assessment" near the top of the document.
# ---------------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
    # bookmark not present - nothing to do
}

# ---------------------------------------------------------------------------
# 2. Merge the three runs that make up question 4 ("Based on the GalaxSee
#    lab...") into a single run. Replacing the full sentence with itself via
#    Find/Replace makes Word collapse the matched range into one run.
# ---------------------------------------------------------------------------
$q4 = "Based on the GalaxSee lab, provide an explanation of what you learned about " + `
      "parallelism and its uses. What results did you expect to see before you began? " + `
      "Did the results that you got confirm your thoughts or were the results different " + `
      "from what you expected? Why do you think this is the case?"
$d.Content.Find.Execute($q4, $true, $false, $false, $false, $false, $true, 1, $false, $q4, 2) | Out-Null

# ---------------------------------------------------------------------------
# Helper: force a run boundary at a given absolute character position by
# toggling Bold on/off across [startPos, startPos+len). This leaves the
# visible formatting untouched, but keeps the engine from silently
# re-merging that span with its identically-formatted neighbors the next
# time some other part of the paragraph is edited.
# ---------------------------------------------------------------------------
function Force-Boundary($startPos, $len) {
    if ($len -le 0) { return }
    $r = $d.Range($startPos, $startPos + $len)
    $r.Bold = $true
    $r.Bold = $false
}

# ---------------------------------------------------------------------------
# 3. "Browse and search the full collection at" -> "...full curriculum at",
#    with "curriculum" split into its own run (same formatting, just its own
#    run - mirrors how Google Docs/Word round-trips a word-level edit).
# ---------------------------------------------------------------------------
$para19 = $d.Paragraphs(19).Range
$p19Start = $para19.Start
$oldWord = "collection"
$newWord = "curriculum"
$idx19 = $para19.Text.IndexOf($oldWord)
$wordStart19 = $p19Start + $idx19
$d.Range($wordStart19, $wordStart19 + $oldWord.Length).Text = $newWord

$prefixLen19 = $idx19
$suffix19 = " at"
Force-Boundary $p19Start $prefixLen19
Force-Boundary $wordStart19 $newWord.Length
# also keep " at" separate from the following <w:br/> run that follows it
Force-Boundary ($wordStart19 + $newWord.Length) $suffix19.Length

# ---------------------------------------------------------------------------
# 4. "material and the rest of the collection in our GitHub repository at"
#    -> "...the curriculum in our GitHub repository at" with "curriculum"
#    split into its own run, and the "_GoBack" bookmark re-inserted
#    immediately after it. This paragraph already had two other run breaks
#    ("We welcome your improvements" | "! You can submit your proposed
#    changes to this " | rest) that must be restored after the text edit,
#    since any in-place text mutation re-flattens identically formatted
#    runs across the whole paragraph.
# ---------------------------------------------------------------------------
$para21 = $d.Paragraphs(21).Range
$p21Start = $para21.Start
$idx21 = $para21.Text.IndexOf($oldWord)
$wordStart21 = $p21Start + $idx21
$d.Range($wordStart21, $wordStart21 + $oldWord.Length).Text = $newWord

$run1Len = "We welcome your improvements".Length
$run2Len = "! You can submit your proposed changes to this ".Length
$run3PrefixLen = $idx21 - $run1Len - $run2Len   # "material and the rest of the "
$suffix21 = " in our GitHub repository at"

Force-Boundary $p21Start $run1Len
Force-Boundary ($p21Start + $run1Len) $run2Len
Force-Boundary ($p21Start + $run1Len + $run2Len) $run3PrefixLen
Force-Boundary $wordStart21 $newWord.Length

$afterCurriculum21 = $wordStart21 + $newWord.Length
$bmRange = $d.Range($afterCurriculum21, $afterCurriculum21)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# keep the trailing " in our GitHub repository at" separate from the following
# <w:br/> run that follows it
Force-Boundary $afterCurriculum21 $suffix21.Length
